$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (AssociateID N1074): update associate name and email
$ws.Range("B6").Value = "Anusha Kodi"
$ws.Range("E6").Value = "anusha.kodi@senecaglobal.com"
